$d = $word.ActiveDocument

# The commit changes a single word in the write-up: "countries" -> "regions"
# within the sentence about the select input on the addPoiForm.php page.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "countries"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "regions"
$find.Forward = $true
$find.Wrap = 0
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $true
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute([ref]$find.Text, $true, $true, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null
